$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# New row 6 - PushP memory usage test on Laptop/Release/Factory_Class branch,
# 116,791 bytes across all three heap-memory test columns.
$ws.Range("A6").Value = 43411.020833333336
$ws.Range("A6").NumberFormat = "m/d/yy h:mm"

$ws.Range("B6").Value = "Laptop"
$ws.Range("C6").Value = "Release"
$ws.Range("D6").Value = "Factory_Class"

$ws.Range("F6").Value = 116791
$ws.Range("F6").NumberFormat = "#,##0"
$ws.Range("G6").Value = 116791
$ws.Range("G6").NumberFormat = "#,##0"
$ws.Range("H6").Value = 116791
$ws.Range("H6").NumberFormat = "#,##0"

$ws.Range("I6").Value = "Implemented CodeList Heap Manager"

# Move the view / selection to match the saved workbook state.
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("I10").Select()
